# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (a per-quarter fund-holding detail sheet,
# matching the layout already used by "2021-Q4") right after "2021-Q4", and
# adds a new leading row to the "总计" (totals) summary sheet for 2022-Q1,
# pushing the existing 2021-Q4 totals row down.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

# --- locate the existing sheets -------------------------------------------------
$q4Sheet = $wb.Worksheets.Item(1)      # "2021-Q4"

# --- 1. create the new "2022-Q1" sheet, right after "2021-Q4" -------------------
$q1Sheet = $wb.Worksheets.Add($null, $q4Sheet)
$q1Sheet.Name = "2022-Q1"

# Fetch the "总计" sheet reference AFTER inserting the new sheet, since its
# position shifted from 2 to 3.
$totalSheet = $wb.Worksheets.Item(3)   # "总计"

# Copy the header-row formatting (bold/border/centered) and the first-column
# index formatting from the existing "2021-Q4" sheet so the new sheet matches
# the established look.
$q4Sheet.Range("B1:H1").Copy()
$q1Sheet.Range("B1:H1").PasteSpecial($xlPasteFormats)

$q4Sheet.Range("A2:A4").Copy()
$q1Sheet.Range("A2:A4").PasteSpecial($xlPasteFormats)

# Header row
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Columns B-G hold text-like values (fund codes/names/figures kept as strings,
# same as the existing "2021-Q4" sheet). Force text entry by pre-formatting as
# "@", then strip the number-format back off afterwards so the cells end up
# with no explicit style (matching the unstyled data cells used elsewhere).
$q1Sheet.Range("B2:G4").NumberFormat = "@"

# Row 2
$q1Sheet.Range("A2").Value = 0
$q1Sheet.Range("B2").Value = "008928"
$q1Sheet.Range("C2").Value = "泰达宏利中证主要消费红利指数A"
$q1Sheet.Range("D2").Value = "5.92"
$q1Sheet.Range("E2").Value = "93.93"
$q1Sheet.Range("F2").Value = "4.56"
$q1Sheet.Range("G2").Value = "0.2700"
$q1Sheet.Range("H2").Value = 8

# Row 3
$q1Sheet.Range("A3").Value = 1
$q1Sheet.Range("B3").Value = "008929"
$q1Sheet.Range("C3").Value = "泰达宏利中证主要消费红利指数C"
$q1Sheet.Range("D3").Value = "3.10"
$q1Sheet.Range("E3").Value = "93.93"
$q1Sheet.Range("F3").Value = "4.56"
$q1Sheet.Range("G3").Value = "0.1414"
$q1Sheet.Range("H3").Value = 8

# Row 4
$q1Sheet.Range("A4").Value = 2
$q1Sheet.Range("B4").Value = "501089"
$q1Sheet.Range("C4").Value = "方正富邦中证主要消费红利指数增强（LOF）"
$q1Sheet.Range("D4").Value = "0.16"
$q1Sheet.Range("E4").Value = "94.62"
$q1Sheet.Range("F4").Value = "4.47"
$q1Sheet.Range("G4").Value = "0.0072"
$q1Sheet.Range("H4").Value = 8

# Drop the temporary "@" formatting now that the text values are locked in.
$q1Sheet.Range("B2:G4").ClearFormats()

# --- 2. add the new 2022-Q1 row to the "总计" sheet, above 2021-Q4 --------------

# Push the existing 2021-Q4 totals row from row 2 down to row 3, carrying its
# formatting (the bold/border/centered style on column A) along with it.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial($xlPasteFormats)

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 0.28

# New first data row: 2022-Q1 totals
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.42
